$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Index Row" value in B2 from 10 to 1
$ws.Range("B2").Value = 1

# Select cell B2 (reflects the active cell/selection change seen in the saved file)
$ws.Range("B2").Select()

$wb.Save()
